# Updated lam and kappa
# Applies updated numeric calibration results to the "Linear" and
# "NonLinear" sheets (including the abs_epsi_autocorr / sig2_1 array
# text stored in B5 / B10 respectively).

$wb = $excel.ActiveWorkbook

# --- Linear sheet ---------------------------------------------------------
$wsLinear = $wb.Worksheets.Item("Linear")

$wsLinear.Range("B2").Value = 0.001162868506602955
$wsLinear.Range("B3").Value = -0.07560068538530179
$wsLinear.Range("B4").Value = 1.652589380335492
$wsLinear.Range("B5").Value = "[0.9999999999999998, 0.22347329692939924, 0.01777152838846454, -0.036304419347339675, -0.039836115734616824, 0.03236204524098861, 0.19886395502798243, 0.3738436546243984, 0.1958957663280312, 0.020029758164960387, -0.03228536244049097, -0.05257059060574654, 0.023530238925371327, 0.19268981643766983, 0.3521225530738798, 0.17944588474245737, 0.007644971051427396, -0.03810666462108624, -0.03723558285148296, 0.005252795641987083]"

# --- NonLinear sheet -------------------------------------------------------
$wsNonLinear = $wb.Worksheets.Item("NonLinear")

$wsNonLinear.Range("B4").Value = -0.00447146869707964
$wsNonLinear.Range("B5").Value = -0.04811823934953421
$wsNonLinear.Range("B7").Value = 0.07570784929995356
$wsNonLinear.Range("B8").Value = -0.2284056395814199
$wsNonLinear.Range("B10").Value = "[1.0, 0.2223258053093444, 0.018421853301887, -0.034475062695248755, -0.03896577827275292, 0.033538029132558524, 0.19717633501802673, 0.37371788223325003, 0.19379714958399782, 0.020139490082570634, -0.03136876274634749, -0.05177716408580365, 0.02326323088384756, 0.1919199982095485, 0.35129581809235955, 0.17877259163892664, 0.008021730551721557, -0.03757255578758416, -0.03481932737571662, 0.005050340906240087]"
